# Refresh the crypto "Price" (column D) and "Volume(1h)" (column E) figures
# on the symbol list, as published by the scheduled GitHub Actions scraper.
# All cells on this sheet are stored as text (even the numeric-looking
# price/percentage values), so each new value is written as a literal string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value (as literal text, matching the
# original inline-string cell contents so prices/percentages keep their
# exact textual formatting, e.g. "321.87" and "8.47%").
$updates = @{
    'D2' = '321.87'
    'E2' = '8.47%'
    'D3' = '51.50'
    'E3' = '23.10%'
    'D4' = '5.360'
    'E4' = '6.97%'
    'D5' = '0.08106'
    'E5' = '7.80%'
    'D6' = '4.572'
    'E6' = '4.32%'
    'D7' = '1.659'
    'E7' = '5.16%'
    'D8' = '1.107'
    'E8' = '19.42%'
    'D9' = '0.1316'
    'E9' = '10.06%'
    'D10' = '0.1956'
    'E10' = '7.06%'
    'D11' = '0.09613'
    'E11' = '8.37%'
    'D12' = '0.04549'
    'E12' = '11.11%'
    'E13' = '0.10%'
    'D14' = '0.001325'
    'E14' = '3.24%'
    'D15' = '0.005850'
    'E15' = '-2.69%'
    'D16' = '3.379'
    'E16' = '0.53%'
    'E17' = '1.26%'
    'D18' = '0.3394'
    'E18' = '2.37%'
    'D19' = '8.192'
    'E19' = '1.41%'
    'D22' = '0.04298'
    'E22' = '4.68%'
    'D23' = '0.001305'
    'E23' = '3.05%'
    'D24' = '0.004267'
    'E24' = '9.75%'
    'D25' = '0.0001348'
    'E25' = '9.47%'
    'D26' = '0.0003715'
    'E26' = '-0.23%'
    'D38' = '0.02780'
    'E38' = '15.53%'
    'D39' = '0.05554'
    'E39' = '6.78%'
    'D40' = '0.006289'
    'E40' = '-0.34%'
    'D41' = '0.007766'
    'E41' = '-0.69%'
    'E42' = '8.96%'
    'D43' = '0.007680'
    'E43' = '3.66%'
    'D44' = '0.008822'
    'E44' = '19.31%'
    'D45' = '0.3523'
    'E45' = '19.44%'
    'D46' = '0.00006829'
    'E46' = '5.65%'
    'E47' = '-0.34%'
    'E48' = '78.17%'
    'D49' = '0.003993'
    'E49' = '-5.08%'
    'D50' = '0.00002096'
    'E50' = '-0.34%'
    'D51' = '0.0001997'
    'E51' = '-0.34%'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text format before assigning so Excel does not silently
    # reinterpret numeric- or percent-looking text as a number (which would
    # turn e.g. "8.47%" into the float 0.0847 instead of keeping the text).
    $cell.NumberFormat = "@"
    $cell.Value2 = $updates[$ref]
}
